# "Minor changes, i'm done for today"
# Refresh the pasted per-member summary stats on the "Per Member Data" sheet
# with the latest computed averages (simple + weighted tele/auto/penalty/total).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Per Member Data")

# row 2 (bredan)
$ws.Range("A2").Value = 70.64285714285714
$ws.Range("B2").Value = 35.333333333333336
$ws.Range("C2").Value = 22.666666666666668
$ws.Range("D2").Value = 0.0
$ws.Range("E2").Value = 70.56109725687725
$ws.Range("F2").Value = 34.60317460335194
$ws.Range("G2").Value = 22.757936507914348
$ws.Range("H2").Value = 0.0

# row 3 (erin)
$ws.Range("A3").Value = 82.16666666666667
$ws.Range("B3").Value = 46.6
$ws.Range("C3").Value = 34.2
$ws.Range("D3").Value = 0.0
$ws.Range("E3").Value = 81.76704545463811
$ws.Range("F3").Value = 46.8911564625177
$ws.Range("G3").Value = 33.44897959201044
$ws.Range("H3").Value = 0.0

# row 4 (luca)
$ws.Range("A4").Value = 69.6
$ws.Range("B4").Value = 26.5
$ws.Range("C4").Value = 29.0
$ws.Range("D4").Value = 0.0
$ws.Range("E4").Value = 68.96774193562257
$ws.Range("F4").Value = 26.5
$ws.Range("G4").Value = 29.0
$ws.Range("H4").Value = 0.0

# row 5 (mason)
$ws.Range("A5").Value = 61.0
$ws.Range("B5").Value = 38.3
$ws.Range("C5").Value = 22.4
$ws.Range("D5").Value = 0.0
$ws.Range("E5").Value = 60.57954545464293
$ws.Range("F5").Value = 37.83219178093086
$ws.Range("G5").Value = 22.571917808179144
$ws.Range("H5").Value = 0.0

# row 6 (zoe)
$ws.Range("A6").Value = 92.75
$ws.Range("B6").Value = 47.333333333333336
$ws.Range("C6").Value = 41.0
$ws.Range("D6").Value = 0.0
$ws.Range("E6").Value = 93.22368421041327
$ws.Range("F6").Value = 48.024691357850564
$ws.Range("G6").Value = 41.148148148110835
$ws.Range("H6").Value = 0.0

# row 7 (cyrus)
$ws.Range("A7").Value = 71.2
$ws.Range("B7").Value = 26.333333333333332
$ws.Range("C7").Value = 28.666666666666668
$ws.Range("D7").Value = 0.0
$ws.Range("E7").Value = 71.03289473687948
$ws.Range("F7").Value = 26.35869565216829
$ws.Range("G7").Value = 28.71739130433658
$ws.Range("H7").Value = 0.0

# row 8 (caleb)
$ws.Range("A8").Value = 61.0
$ws.Range("B8").Value = 38.166666666666664
$ws.Range("C8").Value = 20.166666666666668
$ws.Range("D8").Value = 0.0
$ws.Range("E8").Value = 61.0
$ws.Range("F8").Value = 38.16666666666667
$ws.Range("G8").Value = 20.166666666666664
$ws.Range("H8").Value = 0.0

# row 9 (matt)
$ws.Range("A9").Value = 88.33333333333333
$ws.Range("B9").Value = 47.333333333333336
$ws.Range("C9").Value = 41.0
$ws.Range("D9").Value = 0.0
$ws.Range("E9").Value = 89.1728395059614
$ws.Range("F9").Value = 48.024691357850564
$ws.Range("G9").Value = 41.148148148110835
$ws.Range("H9").Value = 0.0

# row 10 (zach)
$ws.Range("A10").Value = 74.9375
$ws.Range("B10").Value = 33.285714285714285
$ws.Range("C10").Value = 27.0
$ws.Range("D10").Value = 0.0
$ws.Range("E10").Value = 74.01414141434437
$ws.Range("F10").Value = 32.973684210591465
$ws.Range("G10").Value = 26.69736842111582
$ws.Range("H10").Value = 0.0
